$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$ws.Range("L2").Value  = 0.93
$ws.Range("L3").Value  = 1.18
$ws.Range("L4").Value  = 0.95
$ws.Range("L5").Value  = 1.02
$ws.Range("L6").Value  = 0.9
$ws.Range("L7").Value  = 0.96
$ws.Range("L8").Value  = 0.98
$ws.Range("L9").Value  = 0.83
$ws.Range("L10").Value = 0.98
$ws.Range("L12").Value = 1
$ws.Range("L13").Value = 1.16
$ws.Range("L14").Value = 0.9
$ws.Range("L15").Value = 1.16
$ws.Range("L16").Value = 1.04
$ws.Range("L17").Value = 1.15
